$p = $ppt.ActivePresentation
$s = $p.Slides.Item(12)
$shp = $s.Shapes.Item(4)

# Give the body placeholder an explicit position/size (was inheriting from layout).
$shp.Left   = 458787 / 12700
$shp.Top    = 107.375045
$shp.Width  = 8321040 / 12700
$shp.Height = 4935537 / 12700

$tf = $shp.TextFrame
$tr = $tf.TextRange

# Apply text tweaks from the end of the text frame backwards so earlier
# character offsets stay valid while later ones are edited first.
$tr.Characters(344, 30).Text = "(usually smaller and faster)"
$tr.Characters(336, 8).Text  = " XOR EAX "
$tr.Characters(278, 24).Text = "(replace division by 2"
$tr.Characters(271, 7).Text  = " x >> 3  "
$tr.Characters(266, 4).Text  = "x/8        "
$tr.Characters(220, 45).Text = "(replace multiplication by 2 with addition)"
$tr.Characters(214, 6).Text  = " i + i   "
$tr.Characters(171, 3).Text  = " i   "
$tr.Characters(156, 10).Text = "i = i + 1  "
